$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that sits right after
# "introducir" (before " un Id de reserva especifico.").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: extend "El sistema informa al propietario de la anulación."
# with the extra clause about the saved contact method, wrapped with
# proofErr markers (as Word's grammar checker would emit).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("El sistema informa al propietario de la anulación.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $findRng.Start
    $e = $findRng.End
    $replaceRng = $d.Range($s, $e)
    $xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">El sistema </w:t></w:r>
<w:r w:rsidR="002D130F"><w:t>informa al propietario de la anulaci&#243;n</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">por medio del m&#233;todo de contacto </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>guardado</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $replaceRng.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# Change 3: delete the whole "18.a" flow (three paragraphs): the "18.a"
# heading, "Se agota el tiempo para informar al propietario.", and
# "Se marca la tarea como pendiente y se avanza al paso 19."
# ---------------------------------------------------------------------------
$rng18 = $d.Content
$found18 = $rng18.Find.Execute("18.a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found18) {
    $s18 = $rng18.Start
    $p1 = $null
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Start -eq $s18) { $p1 = $para; break }
    }
    if ($p1 -ne $null) {
        $p2 = $p1.Next()
        $p3 = $p2.Next()
        $delRng = $d.Range($p1.Range.Start, $p3.Range.End)
        $delRng.Delete()
    }
}

# ---------------------------------------------------------------------------
# Change 4: delete the "Entre el paso 17 y el paso 18..." paragraph, and
# move the "_GoBack" bookmark to the very start of the following paragraph
# ("Entre el paso 19...").
# ---------------------------------------------------------------------------
$rng17 = $d.Content
$found17 = $rng17.Find.Execute("Entre el paso 17 y el paso 18", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found17) {
    $s17 = $rng17.Start
    $p17 = $null
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Start -eq $s17) { $p17 = $para; break }
    }
    if ($p17 -ne $null) {
        $p17.Range.Delete()
        $afterStart = $p17.Range.Start
        $bmRng = $d.Range($afterStart, $afterStart)
        $d.Bookmarks.Add("_GoBack", $bmRng)
    }
}
